$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at the top of the data block (row 383), shifting the
# existing rows (383-418) down to (387-422), matching the new dimension
# A1:R422.
$ws.Range("A383:R386").EntireRow.Insert()

# New weekly data block for date 45106 (2023-06-29), one row per quality
# grade, inserted at the newly opened rows 383-386.
$newRows = @(
    @{ Row=383; I="Especial"; J=230; K=14000; L=14000; M=14000; P=778 },
    @{ Row=384; I="Primera";  J=740; K=12000; L=13000; M=12351; P=686 },
    @{ Row=385; I="Segunda";  J=370; K=9000;  L=9000;  M=9000;  P=500 },
    @{ Row=386; I="Tercera";  J=100; K=7000;  L=7000;  M=7000;  P=389 }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = 6
    $ws.Cells.Item($r, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($r, 3).Value = "Metropolitana"
    $ws.Cells.Item($r, 4).Value = 45106
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = 100112043
    $ws.Cells.Item($r, 7).Value = "Pepino dulce"
    $ws.Cells.Item($r, 8).Value = "Cultivar IV Región"
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = '$/bandeja 18 kilos'
    $ws.Cells.Item($r, 15).Value = "Provincia de Limarí"
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = 18
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
